# Updated cryptos list with GitHub Actions
# Refresh Price (D) and Volume(1h) (E) columns for the cryptos sheet.
# Rows 20/21 and 39/40 also swap rank position (Chainlink/PEPE, PancakeSwap/Bittensor).
#
# Numeric-looking price strings are written with a leading apostrophe so Excel
# keeps them as literal text (matching the source data, which stores prices as
# text, including values like "1.00" / "0.0000210" whose trailing zeros would
# otherwise be lost if Excel parsed them as numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "90.456.64"
$ws.Range("E2").Value = "  +0.66%  "

$ws.Range("D3").Value = "3.099.35"
$ws.Range("E3").Value = "  +1.33%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "'237.74"
$ws.Range("E5").Value = "  +10.36%  "

$ws.Range("D6").Value = "'624.25"
$ws.Range("E6").Value = "  +2.40%  "

$ws.Range("D7").Value = "'1.13"
$ws.Range("E7").Value = "  +6.61%  "

$ws.Range("D8").Value = "'0.369"
$ws.Range("E8").Value = "  +6.40%  "

$ws.Range("E9").Value = "  +0.01%  "

$ws.Range("D10").Value = "3.100.09"
$ws.Range("E10").Value = "  +1.41%  "

$ws.Range("E11").Value = "  +2.77%  "

$ws.Range("D12").Value = "'0.202"
$ws.Range("E12").Value = "  +3.79%  "

$ws.Range("D13").Value = "'0.0000248"
$ws.Range("E13").Value = "  +4.91%  "

$ws.Range("D14").Value = "'35.02"
$ws.Range("E14").Value = "  +3.25%  "

$ws.Range("D15").Value = "'5.45"
$ws.Range("E15").Value = "  -0.07%  "

$ws.Range("D16").Value = "90.371.26"
$ws.Range("E16").Value = "  +0.79%  "

$ws.Range("D17").Value = "3.685.83"
$ws.Range("E17").Value = "  +1.54%  "

$ws.Range("D18").Value = "3.112.25"
$ws.Range("E18").Value = "  +1.95%  "

$ws.Range("D19").Value = "'3.82"
$ws.Range("E19").Value = "  +2.93%  "

# Row 20 / 21: PEPE and Chainlink swap positions
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "'14.22"
$ws.Range("E20").Value = "  +1.17%  "

$ws.Range("B21").Value = "PEPE"
$ws.Range("C21").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D21").Value = "'0.0000210"
$ws.Range("E21").Value = "  +4.97%  "

$ws.Range("E22").Value = "  +6.52%  "

$ws.Range("D23").Value = "'445.17"
$ws.Range("E23").Value = "  +0.05%  "

$ws.Range("D24").Value = "'8.98"
$ws.Range("E24").Value = "  +2.50%  "

$ws.Range("D25").Value = "'5.88"
$ws.Range("E25").Value = "  +2.13%  "

$ws.Range("D26").Value = "'90.24"
$ws.Range("E26").Value = "  -0.20%  "

$ws.Range("D27").Value = "'12.00"
$ws.Range("E27").Value = "  +1.97%  "

$ws.Range("D28").Value = "3.261.09"
$ws.Range("E28").Value = "  +0.96%  "

$ws.Range("E29").Value = "  +0.08%  "

$ws.Range("D30").Value = "'0.178"
$ws.Range("E30").Value = "  +12.44%  "

$ws.Range("D31").Value = "'0.220"
$ws.Range("E31").Value = "  +10.18%  "

$ws.Range("D32").Value = "'9.15"
$ws.Range("E32").Value = "  +0.83%  "

$ws.Range("D33").Value = "'0.997"
$ws.Range("E33").Value = "  -0.31%  "

$ws.Range("D34").Value = "'0.108"
$ws.Range("E34").Value = "  +28.35%  "

$ws.Range("D35").Value = "'26.31"
$ws.Range("E35").Value = "  -5.98%  "

$ws.Range("E36").Value = "  +39.83%  "

$ws.Range("E37").Value = "  +6.38%  "

$ws.Range("D38").Value = "'7.38"
$ws.Range("E38").Value = "  +10.16%  "

# Row 39 / 40: Bittensor and PancakeSwap swap positions
$ws.Range("B39").Value = "PancakeSwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D39").Value = "'1.91"
$ws.Range("E39").Value = "  +2.34%  "

$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").Value = "'492.49"
$ws.Range("E40").Value = "  -0.02%  "

$ws.Range("D41").Value = "'3.57"
$ws.Range("E41").Value = "  +5.66%  "

$ws.Range("D42").Value = "'1.28"
$ws.Range("E42").Value = "  +1.46%  "

$ws.Range("D43").Value = "'0.416"
$ws.Range("E43").Value = "  -1.16%  "

$ws.Range("D44").Value = "'22.10"
$ws.Range("E44").Value = "  -0.33%  "

$ws.Range("E45").Value = "  +0.03%  "

$ws.Range("D46").Value = "'158.59"
$ws.Range("E46").Value = "  +7.60%  "

$ws.Range("D47").Value = "'1.90"
$ws.Range("E47").Value = "  -0.91%  "

$ws.Range("D48").Value = "'0.676"
$ws.Range("E48").Value = "  -0.21%  "

$ws.Range("D49").Value = "'4.54"
$ws.Range("E49").Value = "  +0.64%  "

$ws.Range("D50").Value = "'44.78"
$ws.Range("E50").Value = "  +0.39%  "

$ws.Range("D51").Value = "'1.32"
$ws.Range("E51").Value = "  +1.17%  "
